{"js": "// Replace the 15 lattice-multiplication problems in the 5x3 table with a new\n// set of problems, keeping the existing per-cell layout:\n//   \"A x B\"\n//   \"  b0    b1\"   (digits of B, spaced)\n//   \"  ----\"\n//   \"a0|    |\"     (tens digit of A)\n//   \"a1|    |\"     (ones digit of A)\n\nconst newProblems = [\n  [75, 81], [80, 37], [89, 86],\n  [28, 53], [60, 17], [30, 91],\n  [67, 27], [58, 61], [64, 96],\n  [66, 12], [30, 67], [95, 66],\n  [50, 69], [79, 20], [59, 59],\n];\n\nfunction buildCellOoxml(a, b) {\n  const aStr = String(a).padStart(2, \"0\");\n  const bStr = String(b).padStart(2, \"0\");\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r>' +\n    '<w:rPr><w:sz w:val=\"32\"/></w:rPr>' +\n    \"<w:t>\" + a + \" x \" + b + \"</w:t>\" +\n    \"<w:br/>\" +\n    '<w:t xml:space=\"preserve\">  ' + bStr[0] + \"    \" + bStr[1] + \"</w:t>\" +\n    \"<w:br/>\" +\n    '<w:t xml:space=\"preserve\">  ----</w:t>' +\n    \"<w:br/>\" +\n    \"<w:t>\" + aStr[0] + \"|    |</w:t>\" +\n    \"<w:br/>\" +\n    \"<w:t>\" + aStr[1] + \"|    |</w:t>\" +\n    \"</w:r>\" +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst colCount = 3;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const [a, b] = newProblems[idx];\n    const cell = table.getCell(r, c);\n    cell.body.insertOoxml(buildCellOoxml(a, b), Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems in the 5x3 table with a new\n# set of problems, keeping the existing per-cell layout:\n#   \"A x B\"\n#   \"  b0    b1\"   (digits of B, spaced)\n#   \"  ----\"\n#   \"a0|    |\"     (tens digit of A)\n#   \"a1|    |\"     (ones digit of A)\n#\n# NOTE: string interpolation (\"$a$b\") is used instead of the \"+\" operator to\n# join text fragments, because this host's PowerShell engine silently adds\n# two numeric-looking strings (e.g. \"  8\" + \"    1\") as numbers instead of\n# concatenating them.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$problems = @(\n    @(75,81), @(80,37), @(89,86),\n    @(28,53), @(60,17), @(30,91),\n    @(67,27), @(58,61), @(64,96),\n    @(66,12), @(30,67), @(95,66),\n    @(50,69), @(79,20), @(59,59)\n)\n\n$lineBreak = [char]11\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $problems[$idx]\n        $a = $pair[0]\n        $b = $pair[1]\n        $aStr = \"{0:D2}\" -f $a\n        $bStr = \"{0:D2}\" -f $b\n        $aTens = $aStr.Substring(0,1)\n        $aOnes = $aStr.Substring(1,1)\n        $bTens = $bStr.Substring(0,1)\n        $bOnes = $bStr.Substring(1,1)\n\n        $line1 = \"$a x $b\"\n        $line2 = \"  $bTens    $bOnes\"\n        $line3 = \"  ----\"\n        $line4 = \"$aTens|    |\"\n        $line5 = \"$aOnes|    |\"\n\n        $newText = \"$line1$lineBreak$line2$lineBreak$line3$lineBreak$line4$lineBreak$line5\"\n\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newText\n\n        $idx++\n    }\n}\n"}
